# Trade #28 closed at 2026-02-17 23:58:06 - unknown UNKNOWN +0.000%
#
# Updates the rolling performance summary after a new (28th) trade closed,
# and appends that trade's row to the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - roll the aggregate stats forward
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1501.1    # Current Capital
$summary.Range("B4").Value = 1.1       # Total P&L $
$summary.Range("B5").Value = 0.79      # Total P&L %
$summary.Range("B6").Value = 28        # Total Trades
$summary.Range("B8").Value = 10        # Losing Trades
$summary.Range("B9").Value = 57.14     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 101.1      # Capital
$status.Range("D6").Value = 28         # Trades
$status.Range("E6").Value = 1.1        # P&L $
$status.Range("F6").Value = 1.1        # P&L %
$status.Range("G6").Value = 57.14      # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new trade row (row 29) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
$newRow = @{
    A = 28
    B = "2026-02-17"
    C = "23:57:59"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.5600000000000001
    G = 0.53
    H = "CLOSED"
    I = -5.3571
    J = -0.03
    K = 101.1
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(29, 1).Value = $newRow.A

    # Date/time-shaped text needs to be forced to Text format first,
    # otherwise Excel auto-parses "2026-02-17" into a date serial.
    $ws.Cells.Item(29, 2).NumberFormat = "@"
    $ws.Cells.Item(29, 2).Value = $newRow.B
    $ws.Cells.Item(29, 2).ClearFormats()

    $ws.Cells.Item(29, 3).Value = $newRow.C
    $ws.Cells.Item(29, 4).Value = $newRow.D
    $ws.Cells.Item(29, 5).Value = $newRow.E
    $ws.Cells.Item(29, 6).Value = $newRow.F
    $ws.Cells.Item(29, 7).Value = $newRow.G
    $ws.Cells.Item(29, 8).Value = $newRow.H
    $ws.Cells.Item(29, 9).Value = $newRow.I
    $ws.Cells.Item(29, 10).Value = $newRow.J
    $ws.Cells.Item(29, 11).Value = $newRow.K
    $ws.Cells.Item(29, 12).Value = $newRow.L
    $ws.Cells.Item(29, 13).Value = $newRow.M
    $ws.Cells.Item(29, 14).Value = $newRow.N
    $ws.Cells.Item(29, 15).Value = $newRow.O
    $ws.Cells.Item(29, 16).Value = $newRow.P
    $ws.Cells.Item(29, 17).Value = $newRow.Q
}
